$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the duplicate "1-5 scale" Units values in column N (rows 2-4).
# These strings ("Predation 1-5 scale", "CrBr 1-5 scale", "CrBuCt 1-5 scale")
# become unused shared strings and are dropped from the workbook on save.
$ws.Range("N2:N4").ClearContents()

# Update the view: scroll so column C is the left-most visible column, and
# select N2:N4 with N2 as the active cell.
$ws.Application.ActiveWindow.ScrollColumn = 3
$ws.Range("N2:N4").Select()
